$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the SQL text in B7: the ORDER BY clause now sorts on srv.survival_id
# instead of prt.participant_id.
$old = $ws.Range("B7").Value2
$new = $old -replace "ORDER BY \r?\n    prt\.participant_id ASC", "ORDER BY `n    srv.survival_id ASC"
$ws.Range("B7").Value = $new

# Nudge the font size back to the workbook default (11pt) for this cell while
# keeping the existing wrap-text alignment.
$ws.Range("B7").Font.Size = 11
$ws.Range("B7").Font.ThemeColor = 1

# Scroll the view back to show column A.
$ws.Application.ActiveWindow.ScrollColumn = 1
